$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (avoids numeric/date auto-coercion
# for number-like strings such as "65.037.52" or "0.493"), then drop the temporary
# "@" number-format style so the cell stays style-less like the original.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '65.037.52'
Set-TextValue 'E2' '  -0.28%  '

Set-TextValue 'D3' '3.540.58'
Set-TextValue 'E3' '  +0.25%  '

Set-TextValue 'E4' '  +0.01%  '

Set-TextValue 'D5' '598.04'
Set-TextValue 'E5' '  +0.51%  '

Set-TextValue 'D6' '134.33'
Set-TextValue 'E6' '  -3.42%  '

Set-TextValue 'D7' '3.533.92'
Set-TextValue 'E7' '  -0.01%  '

Set-TextValue 'E8' '  +0.04%  '

Set-TextValue 'D9' '0.493'
Set-TextValue 'E9' '  -0.39%  '

Set-TextValue 'D10' '0.123'
Set-TextValue 'E10' '  -2.13%  '

Set-TextValue 'D11' '6.97'
Set-TextValue 'E11' '  -3.32%  '

Set-TextValue 'D12' '0.386'
Set-TextValue 'E12' '  -0.92%  '

Set-TextValue 'D13' '4.143.37'
Set-TextValue 'E13' '  +0.34%  '

Set-TextValue 'D14' '0.0000182'
Set-TextValue 'E14' '  -2.74%  '

Set-TextValue 'D15' '26.93'
Set-TextValue 'E15' '  -0.54%  '

Set-TextValue 'D16' '3.546.82'
Set-TextValue 'E16' '  +0.42%  '

Set-TextValue 'E17' '  +0.12%  '

Set-TextValue 'D18' '65.168.61'
Set-TextValue 'E18' '  +0.17%  '

Set-TextValue 'D19' '9.92'
Set-TextValue 'E19' '  -2.04%  '

Set-TextValue 'D20' '14.36'
Set-TextValue 'E20' '  +0.66%  '

Set-TextValue 'D21' '5.83'
Set-TextValue 'E21' '  -0.73%  '

Set-TextValue 'D22' '390.92'
Set-TextValue 'E22' '  -1.20%  '

Set-TextValue 'D23' '0.577'
Set-TextValue 'E23' '  +0.79%  '

Set-TextValue 'D24' '3.685.33'
Set-TextValue 'E24' '  +0.30%  '

Set-TextValue 'D25' '73.90'
Set-TextValue 'E25' '  -0.46%  '

Set-TextValue 'E26' '  +0.14%  '

Set-TextValue 'D27' '0.0000112'
Set-TextValue 'E27' '  -2.56%  '

Set-TextValue 'D28' '7.68'
Set-TextValue 'E28' '  -1.37%  '

Set-TextValue 'B29' 'Binance-PegBSC-USD'
Set-TextValue 'C29' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D29' '1.00'
Set-TextValue 'E29' '  +0.24%  '

Set-TextValue 'B30' 'Fetch.AI'
Set-TextValue 'C30' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D30' '1.54'
Set-TextValue 'E30' '  +25.04%  '

Set-TextValue 'B31' 'InternetComputer(DFINITY)'
Set-TextValue 'C31' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D31' '8.48'
Set-TextValue 'E31' '  +2.19%  '

Set-TextValue 'D32' '2.28'
Set-TextValue 'E32' '  +1.22%  '

Set-TextValue 'D33' '3.542.79'
Set-TextValue 'E33' '  -0.18%  '

Set-TextValue 'D34' '23.99'
Set-TextValue 'E34' '  +0.71%  '

Set-TextValue 'E35' '  +0.00%  '

Set-TextValue 'D36' '0.145'
Set-TextValue 'E36' '  +0.39%  '

Set-TextValue 'D37' '6.94'
Set-TextValue 'E37' '  -0.56%  '

Set-TextValue 'B38' 'Monero'
Set-TextValue 'C38' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D38' '169.01'
Set-TextValue 'E38' '  -0.39%  '

Set-TextValue 'B39' 'ImmutableX'
Set-TextValue 'C39' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D39' '1.55'
Set-TextValue 'E39' '  +0.96%  '

Set-TextValue 'D40' '5.00'
Set-TextValue 'E40' '  +2.36%  '

Set-TextValue 'D41' '0.0805'
Set-TextValue 'E41' '  +0.70%  '

Set-TextValue 'D42' '0.824'
Set-TextValue 'E42' '  +0.30%  '

Set-TextValue 'D43' '25.96'
Set-TextValue 'E43' '  -2.78%  '

Set-TextValue 'D44' '42.95'
Set-TextValue 'E44' '  +0.00%  '

Set-TextValue 'E45' '  +3.76%  '

Set-TextValue 'E46' '  -0.18%  '

Set-TextValue 'D47' '4.44'
Set-TextValue 'E47' '  +0.28%  '

Set-TextValue 'D48' '1.65'
Set-TextValue 'E48' '  -0.32%  '

Set-TextValue 'D49' '2.460.74'
Set-TextValue 'E49' '  +5.50%  '

Set-TextValue 'D50' '6.91'
Set-TextValue 'E50' '  +1.67%  '

Set-TextValue 'D51' '0.0263'
Set-TextValue 'E51' '  +1.76%  '
